# Add files via upload
# - Adds a new daily snapshot column (CP) to Sheet1, dated "15-nov" (the
#   day after the existing last column "14-nov"), copying the latest
#   VLOOKUP-derived values (same values now found in CB/CC).
# - Updates the underlying raw lookup table on Sheet3 ($A$20:$B$36), which
#   drives every VLOOKUP formula in Sheet1 (CB/CC) and Sheet3 (C2:C18).
# - Moves the saved cell selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# 1. Refresh the raw source values on Sheet3 ($A$20:$B$36). Everything
#    else (Sheet1 CB/CC, Sheet3 C2:C18) is a VLOOKUP against this table
#    and recalculates automatically.
# ---------------------------------------------------------------------
$ws3.Range("B20").Value = 11.504357031855415   # 3D QUESO 92GX27
$ws3.Range("B23").Value = 6.5273824224214518   # DORITOS QUESO 70X40G
$ws3.Range("B24").Value = 6.3626170798673742   # DORITOS QUESO 77GX26
$ws3.Range("B26").Value = 7.700471435593756    # LAYS CLASICAS 145GRX18
$ws3.Range("B27").Value = 4.3630321768791882   # LAYS CLASICAS 249GRX14
$ws3.Range("B28").Value = 5.5581799387251509   # LAYS CLASICAS 40GX68
$ws3.Range("B29").Value = 4.8888069949958188   # LAYS CLASICAS 94GRX25
$ws3.Range("B30").Value = 0                    # LAYS ONDAS FH 30GX72
$ws3.Range("B31").Value = 0.6349419999999687   # LAYS ONDAS FH 70GX28
$ws3.Range("B32").Value = 16.172460150894779   # LAYS QSO Y CEBOLLA 34GX72
$ws3.Range("B33").Value = 10.729297437759444   # PEHUAMAR ACANALADA 520GX9
$ws3.Range("B34").Value = 6.6093339900976211   # PEHUAMAR MAICITOS 285GX10
$ws3.Range("B35").Value = 8.9552349995673524   # PEHUAMAR PAPA LISA 520GX9
$ws3.Range("B36").Value = 44.399509122518197   # QUAKER AVENA INSTANT FORTIF 18X280G

# ---------------------------------------------------------------------
# 2. Append the new "15-nov" snapshot column (CP) on Sheet1, mirroring
#    the existing "14-nov" column (CO): a text header row and, below it,
#    the latest lookup value (now equal to CB/CC) per product row.
# ---------------------------------------------------------------------
$ws1.Range("CP1").Value = "15-nov"
$ws1.Range("CP1").NumberFormat = "@"

$ws1.Range("CP2").NumberFormat = "0"
$ws1.Range("CP3").NumberFormat = "0"
$ws1.Range("CP4").NumberFormat = "0"
$ws1.Range("CP5").NumberFormat = "0"
$ws1.Range("CP6").NumberFormat = "0"
$ws1.Range("CP7").NumberFormat = "0"
$ws1.Range("CP8").NumberFormat = "0"
$ws1.Range("CP9").NumberFormat = "0"
$ws1.Range("CP10").NumberFormat = "0"
$ws1.Range("CP11").NumberFormat = "0"
$ws1.Range("CP12").NumberFormat = "0"
$ws1.Range("CP13").NumberFormat = "0"
$ws1.Range("CP14").NumberFormat = "0"
$ws1.Range("CP15").NumberFormat = "0"
$ws1.Range("CP16").NumberFormat = "0"
$ws1.Range("CP17").NumberFormat = "0"
$ws1.Range("CP18").NumberFormat = "0"

$ws1.Range("CP2").Value = 5.5581799387251509
$ws1.Range("CP3").Value = 4.8888069949958188
$ws1.Range("CP4").Value = 7.700471435593756
$ws1.Range("CP5").Value = 4.3630321768791882
$ws1.Range("CP6").Value = 6.5273824224214518
$ws1.Range("CP7").Value = 6.3626170798673742
$ws1.Range("CP8").Value = 0
$ws1.Range("CP9").Value = 8.9552349995673524
$ws1.Range("CP10").Value = 10.729297437759444
$ws1.Range("CP11").Value = 6.6093339900976211
$ws1.Range("CP12").Value = 11.504357031855415
$ws1.Range("CP13").Value = 0
$ws1.Range("CP14").Value = 44.399509122518197
$ws1.Range("CP15").Value = 16.172460150894779
$ws1.Range("CP16").Value = 0
$ws1.Range("CP17").Value = 0
$ws1.Range("CP18").Value = 0.6349419999999687

# ---------------------------------------------------------------------
# 3. Move the saved selection on Sheet1 to the newly-edited area.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("CL26:CL27").Select()
